$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.178.74'
$ws.Range("E2").Value = '''  +1.99%  '

$ws.Range("D3").Value = '''3.371.26'
$ws.Range("E3").Value = '''  +2.00%  '

$ws.Range("D4").Value = '''0.997'
$ws.Range("E4").Value = '''  -0.29%  '

$ws.Range("D5").Value = '''590.11'
$ws.Range("E5").Value = '''  +6.34%  '

$ws.Range("D6").Value = '''189.09'
$ws.Range("E6").Value = '''  +0.30%  '

$ws.Range("D7").Value = '''0.609'
$ws.Range("E7").Value = '''  +5.04%  '

$ws.Range("E8").Value = '''  +0.00%  '

$ws.Range("D9").Value = '''0.187'
$ws.Range("E9").Value = '''  +3.99%  '

$ws.Range("D10").Value = '''0.591'
$ws.Range("E10").Value = '''  +1.46%  '

$ws.Range("D11").Value = '''47.69'
$ws.Range("E11").Value = '''  +1.85%  '

$ws.Range("D12").Value = '''0.0000275'
$ws.Range("E12").Value = '''  +3.36%  '

$ws.Range("D13").Value = '''655.06'
$ws.Range("E13").Value = '''  +9.50%  '

$ws.Range("D14").Value = '''3.893.48'
$ws.Range("E14").Value = '''  +1.49%  '

$ws.Range("D15").Value = '''8.62'
$ws.Range("E15").Value = '''  +0.11%  '

$ws.Range("D16").Value = '''67.055.73'
$ws.Range("E16").Value = '''  +1.68%  '

$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").Value = '''18.09'
$ws.Range("E17").Value = '''  +1.22%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '''0.119'
$ws.Range("E18").Value = '''  +0.83%  '

$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '''3.358.11'
$ws.Range("E19").Value = '''  +1.38%  '

$ws.Range("D20").Value = '''11.24'
$ws.Range("E20").Value = '''  +2.15%  '

$ws.Range("D21").Value = '''0.910'
$ws.Range("E21").Value = '''  +1.60%  '

$ws.Range("D22").Value = '''18.02'
$ws.Range("E22").Value = '''  -2.50%  '

$ws.Range("D23").Value = '''5.14'
$ws.Range("E23").Value = '''  +1.00%  '

$ws.Range("D24").Value = '''101.47'
$ws.Range("E24").Value = '''  +0.57%  '

$ws.Range("D25").Value = '''4.03'
$ws.Range("E25").Value = '''  +2.25%  '

$ws.Range("D26").Value = '''2.84'
$ws.Range("E26").Value = '''  +3.90%  '

$ws.Range("D27").Value = '''9.79'
$ws.Range("E27").Value = '''  +3.15%  '

$ws.Range("D28").Value = '''32.38'
$ws.Range("E28").Value = '''  +7.37%  '

$ws.Range("D29").Value = '''8.75'
$ws.Range("E29").Value = '''  +1.09%  '

$ws.Range("D30").Value = '''6.95'
$ws.Range("E30").Value = '''  +3.59%  '

$ws.Range("D31").Value = '''616.58'
$ws.Range("E31").Value = '''  +7.71%  '

$ws.Range("D32").Value = '''3.96'
$ws.Range("E32").Value = '''  +3.81%  '

$ws.Range("D33").Value = '''11.26'
$ws.Range("E33").Value = '''  +2.36%  '

$ws.Range("D34").Value = '''3.894.60'
$ws.Range("E34").Value = '''  +5.38%  '

$ws.Range("D35").Value = '''0.107'
$ws.Range("E35").Value = '''  +3.00%  '

$ws.Range("E36").Value = '''  +0.04%  '

$ws.Range("D37").Value = '''55.70'
$ws.Range("E37").Value = '''  -2.16%  '

$ws.Range("D38").Value = '''2.81'
$ws.Range("E38").Value = '''  +5.08%  '

$ws.Range("D39").Value = '''0.132'
$ws.Range("E39").Value = '''  +2.91%  '

$ws.Range("D40").Value = '''33.95'
$ws.Range("E40").Value = '''  +0.81%  '

$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").Value = '''0.0₃0715'
$ws.Range("E41").Value = '''  +2.89%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''3.25'
$ws.Range("E42").Value = '''  +0.12%  '

$ws.Range("D43").Value = '''0.346'
$ws.Range("E43").Value = '''  +2.46%  '

$ws.Range("D44").Value = '''3.39'
$ws.Range("E44").Value = '''  -0.60%  '

$ws.Range("D45").Value = '''0.0425'
$ws.Range("E45").Value = '''  +1.83%  '

$ws.Range("D46").Value = '''0.132'
$ws.Range("E46").Value = '''  +2.35%  '

$ws.Range("D47").Value = '''2.59'
$ws.Range("E47").Value = '''  +1.66%  '

$ws.Range("E48").Value = '''  +0.07%  '

$ws.Range("B49").Value = 'CoreDAO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D49").Value = '''2.90'
$ws.Range("E49").Value = '''  -15.90%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '''1.35'
$ws.Range("E50").Value = '''  +8.77%  '

$ws.Range("D51").Value = '''130.75'
$ws.Range("E51").Value = '''  +6.09%  '
